# Add the new "InvalidLogin" worksheet right after "ValidLogin" and
# populate it with the invalid-login test data, then make it the active tab.

$wb = $excel.ActiveWorkbook

# Create the new sheet and give it its final name.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "InvalidLogin"

# Put it right after the existing "ValidLogin" sheet.
$newSheet.Move($null, $wb.Worksheets.Item("ValidLogin"))

# Re-fetch a fresh reference to the sheet by name (the handle captured
# before the Move can go stale) and use that for every subsequent write.
$inv = $wb.Worksheets.Item("InvalidLogin")

# Header row + first invalid-login test row (Username / Password columns).
$inv.Range("A1").Value = "Username"
$inv.Range("B1").Value = "Password"
$inv.Range("A2").Value = "abcd"
$inv.Range("B2").Value = "xyz"

# FailMsg column (header + the message repeated for every test row).
$inv.Range("C1").Value = "FailMsg"
$inv.Range("C2").Value = "Err Msg is Not Dispalyed"
$inv.Range("C3").Value = "Err Msg is Not Dispalyed"
$inv.Range("C4").Value = "Err Msg is Not Dispalyed"
$inv.Range("C5").Value = "Err Msg is Not Dispalyed"

# Remaining invalid-login test rows.
$inv.Range("A3").Value = "admin"
$inv.Range("B3").Value = "damager"
$inv.Range("A4").Value = "admin"
$inv.Range("B5").Value = "manager"

# Auto-fit the Password / FailMsg columns to their new content.
$inv.Columns.Item(2).AutoFit() | Out-Null
$inv.Columns.Item(3).AutoFit() | Out-Null

# Make the new sheet the active tab, with the last-entered cell selected,
# and zoom in like the source workbook.
$excel.Worksheets("InvalidLogin").Activate()
$inv.Range("C5").Select() | Out-Null
$excel.ActiveWindow.Zoom = 220
